$wb = $excel.ActiveWorkbook

# --- Sheet "Overview": Status text for zh-cn / de-de columns (E2, F2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# Widen the Status columns (E, F) to fit the new, longer text.
$wsOverview.Range("E:E").ColumnWidth = 29.166666666666668
$wsOverview.Range("F:F").ColumnWidth = 29.166666666666668

# --- Sheet "zh-cn": Status text, handback datetime, clear stale error ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-10-10 09:40:25"

# Clear the now-stale "handback not latest" error message, keeping the cell
# as an (empty) text cell rather than removing it outright.
$wsZhCn.Range("P2").Value = "'"
$wsZhCn.Range("P2").Style = "Normal"

# Resize Status column (wider) and Error Detail column (narrower, now empty).
$wsZhCn.Range("C:C").ColumnWidth = 29.166666666666668
$wsZhCn.Range("P:P").ColumnWidth = 12.833333333333334

# --- Sheet "de-de": Status text, handback datetime, clear stale error ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-10-10 09:40:41"

$wsDeDe.Range("P2").Value = "'"
$wsDeDe.Range("P2").Style = "Normal"

$wsDeDe.Range("C:C").ColumnWidth = 29.166666666666668
$wsDeDe.Range("P:P").ColumnWidth = 12.833333333333334
